# Added Indian MF 1st Stab
# This script inserts 9 new "week" columns (the most recent weeks of
# analyst-rating history) in front of the existing data table, fills in
# the new header dates and filler ("UN") values, records a new
# "Set Price Target" event for Morgan Stanley in the newly added
# "Jul_17" week, and normalises the few short rows (30-33) so that their
# filler cells stay contiguous instead of being pushed out by the insert.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert 9 new columns before column B (the first data column).
#    Everything that used to live in B:V now lives in K:AE.
# ---------------------------------------------------------------------
$ws.Range("B1:J1").EntireColumn.Insert()

# ---------------------------------------------------------------------
# 2. New header row (most recent week first, descending).
# ---------------------------------------------------------------------
$ws.Range("B1").Value = "Sep_08"
$ws.Range("C1").Value = "Aug_25"
$ws.Range("D1").Value = "Aug_04"
$ws.Range("E1").Value = "Jul_23"
$ws.Range("F1").Value = "Jul_17"
$ws.Range("G1").Value = "Jul_07"
$ws.Range("H1").Value = "Jun_30"
$ws.Range("I1").Value = "Jun_24"
$ws.Range("J1").Value = "Jun_16"

# ---------------------------------------------------------------------
# 3. Fill the new week columns for every analyst row (2-29) with the
#    "no rating change" placeholder used throughout the sheet.
# ---------------------------------------------------------------------
$ws.Range("B2:J29").Value = "UN"

# ---------------------------------------------------------------------
# 4. Record the new analyst event: Morgan Stanley (row 7) set a price
#    target on 7/15/2019, which falls inside the new "Jul_17" week
#    (column F). Highlight it the same way other positive events
#    (upgrades / price-target boosts) are highlighted elsewhere in the
#    sheet.
# ---------------------------------------------------------------------
$ws.Range("F7").Value = "7/15/2019,Set Price Target,Buy,`$71.00"
$ws.Range("F7").Interior.ColorIndex = 35

# ---------------------------------------------------------------------
# 5. Rows 30-33 only ever had data for a handful of the earliest weeks
#    (they don't extend across the whole table). In the source report
#    these rows are not shifted by the insert - their existing filler
#    cells stay put and the 9 new filler cells are simply appended
#    right after them. Undo the automatic shift for just these rows.
# ---------------------------------------------------------------------
$ws.Range("B30:AE31").ClearContents()
$ws.Range("B30:Y31").Value = "UN"

$ws.Range("B32:AE33").ClearContents()
$ws.Range("B32:P33").Value = "UN"
